$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("G370").Value = 24.793805000000003
$ws.Range("G370").NumberFormat = "0.00"
